# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Valor Mora" (F) column values for the existing "Periodo Mora" rows
# got refreshed from the source database: rows 16/18 swap their amounts,
# and rows 19/22 swap their amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = 15625
$ws.Range("F18").Value = 31249

$ws.Range("F19").Value = 20800
$ws.Range("F22").Value = 52000
